# Auto-generated PowerShell Excel COM-interop script
# Applies numeric value updates to the Hades_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1047.619
$ws.Cells.Item(40, 10).Value = 1047.619
$ws.Cells.Item(40, 12).Value = 1047.619
$ws.Cells.Item(40, 14).Value = -1397.619
$ws.Cells.Item(132, 8).Value = 2165.547
$ws.Cells.Item(132, 9).Value = 1783.0222
$ws.Cells.Item(132, 10).Value = 4317.25
$ws.Cells.Item(132, 11).Value = 5349.0666
$ws.Cells.Item(132, 12).Value = 12951.75
$ws.Cells.Item(132, 13).Value = -2819.0666
$ws.Cells.Item(132, 14).Value = -18011.75
$ws.Cells.Item(137, 8).Value = 1821630.8
$ws.Cells.Item(137, 9).Value = 4001143.5
$ws.Cells.Item(137, 10).Value = 5370.067
$ws.Cells.Item(137, 11).Value = 12003430.5
$ws.Cells.Item(137, 12).Value = 16110.201
$ws.Cells.Item(137, 13).Value = -12000880.5
$ws.Cells.Item(137, 14).Value = -21210.201
$ws.Cells.Item(138, 8).Value = 3847635.5
$ws.Cells.Item(138, 9).Value = 833.06665
$ws.Cells.Item(138, 10).Value = 5407150
$ws.Cells.Item(138, 11).Value = 2499.19995
$ws.Cells.Item(138, 12).Value = 16221450
$ws.Cells.Item(138, 13).Value = 2640.80005
$ws.Cells.Item(138, 14).Value = -16231730
$ws.Cells.Item(141, 8).Value = 3741.4
$ws.Cells.Item(141, 9).Value = 3669.875
$ws.Cells.Item(141, 10).Value = 4027.5
$ws.Cells.Item(141, 11).Value = 11009.625
$ws.Cells.Item(141, 12).Value = 12082.5
$ws.Cells.Item(141, 13).Value = -5829.625
$ws.Cells.Item(141, 14).Value = -22442.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1351.289
$ws.Cells.Item(2, 9).Value = 1220.9474
$ws.Cells.Item(2, 11).Value = 1220.9474
$ws.Cells.Item(2, 13).Value = -1107.9474
$ws.Cells.Item(61, 8).Value = 30365768
$ws.Cells.Item(61, 9).Value = 45501268
$ws.Cells.Item(61, 11).Value = 45501268
$ws.Cells.Item(61, 13).Value = -45501056
$ws.Cells.Item(74, 8).Value = 5782721
$ws.Cells.Item(74, 9).Value = 6966092.5
$ws.Cells.Item(74, 10).Value = 102537.8
$ws.Cells.Item(74, 11).Value = 6966092.5
$ws.Cells.Item(74, 12).Value = 102537.8
$ws.Cells.Item(74, 13).Value = -6965218.5
$ws.Cells.Item(74, 14).Value = -104285.8
$ws.Cells.Item(77, 8).Value = 5782721
$ws.Cells.Item(77, 9).Value = 6966092.5
$ws.Cells.Item(77, 10).Value = 102537.8
$ws.Cells.Item(77, 11).Value = 34830462.5
$ws.Cells.Item(77, 12).Value = 512689
$ws.Cells.Item(77, 13).Value = -34826094.5
$ws.Cells.Item(77, 14).Value = -521425
$ws.Cells.Item(116, 8).Value = 1351.289
$ws.Cells.Item(116, 9).Value = 1220.9474
$ws.Cells.Item(116, 11).Value = 1220.9474
$ws.Cells.Item(116, 13).Value = 1073.0526
$ws.Cells.Item(132, 8).Value = 38959.035
$ws.Cells.Item(132, 9).Value = 24136.477
$ws.Cells.Item(132, 10).Value = 93308.414
$ws.Cells.Item(132, 11).Value = 72409.431
$ws.Cells.Item(132, 12).Value = 279925.242
$ws.Cells.Item(132, 13).Value = -69879.431
$ws.Cells.Item(132, 14).Value = -284985.242
$ws.Cells.Item(136, 8).Value = 30365768
$ws.Cells.Item(136, 9).Value = 45501268
$ws.Cells.Item(136, 11).Value = 136503804
$ws.Cells.Item(136, 13).Value = -136501254

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1351.289
$ws.Cells.Item(3, 9).Value = 1220.9474
$ws.Cells.Item(3, 11).Value = 1220.9474
$ws.Cells.Item(3, 13).Value = -1106.9474
$ws.Cells.Item(64, 8).Value = 1125.125
$ws.Cells.Item(64, 9).Value = 829.5
$ws.Cells.Item(64, 10).Value = 1420.75
$ws.Cells.Item(64, 11).Value = 829.5
$ws.Cells.Item(64, 12).Value = 1420.75
$ws.Cells.Item(64, 13).Value = -604.5
$ws.Cells.Item(64, 14).Value = -1870.75
$ws.Cells.Item(67, 8).Value = 1125.125
$ws.Cells.Item(67, 9).Value = 829.5
$ws.Cells.Item(67, 10).Value = 1420.75
$ws.Cells.Item(67, 11).Value = 829.5
$ws.Cells.Item(67, 12).Value = 1420.75
$ws.Cells.Item(67, 13).Value = -49.5
$ws.Cells.Item(67, 14).Value = -2980.75
$ws.Cells.Item(94, 8).Value = 1037.3846
$ws.Cells.Item(94, 9).Value = 865.0909
$ws.Cells.Item(94, 11).Value = 865.0909
$ws.Cells.Item(94, 13).Value = -414.0909

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 43479900
$ws.Cells.Item(58, 9).Value = 62501720
$ws.Cells.Item(58, 10).Value = 1457.7142
$ws.Cells.Item(58, 11).Value = 62501720
$ws.Cells.Item(58, 12).Value = 1457.7142
$ws.Cells.Item(58, 13).Value = -62501517
$ws.Cells.Item(58, 14).Value = -1863.7142
$ws.Cells.Item(132, 8).Value = 30771.742
$ws.Cells.Item(132, 9).Value = 1853.2693
$ws.Cells.Item(132, 11).Value = 5559.8079
$ws.Cells.Item(132, 13).Value = -3029.8079
$ws.Cells.Item(134, 8).Value = 20366.887
$ws.Cells.Item(134, 9).Value = 1373.1395
$ws.Cells.Item(134, 10).Value = 102040
$ws.Cells.Item(134, 11).Value = 4119.4185
$ws.Cells.Item(134, 12).Value = 306120
$ws.Cells.Item(134, 13).Value = -1584.4185
$ws.Cells.Item(134, 14).Value = -311190
$ws.Cells.Item(136, 8).Value = 43479900
$ws.Cells.Item(136, 9).Value = 62501720
$ws.Cells.Item(136, 10).Value = 1457.7142
$ws.Cells.Item(136, 11).Value = 187505160
$ws.Cells.Item(136, 12).Value = 4373.142599999999
$ws.Cells.Item(136, 13).Value = -187502610
$ws.Cells.Item(136, 14).Value = -9473.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 807.375
$ws.Cells.Item(92, 9).Value = 932.5
$ws.Cells.Item(92, 10).Value = 682.25
$ws.Cells.Item(92, 11).Value = 2797.5
$ws.Cells.Item(92, 12).Value = 2046.75
$ws.Cells.Item(92, 13).Value = -1549.5
$ws.Cells.Item(92, 14).Value = -4542.75
$ws.Cells.Item(95, 8).Value = 6250
$ws.Cells.Item(95, 10).Value = 6250
$ws.Cells.Item(95, 12).Value = 18750
$ws.Cells.Item(95, 14).Value = -22868
$ws.Cells.Item(129, 8).Value = 3625374.2
$ws.Cells.Item(129, 9).Value = 1789.7778
$ws.Cells.Item(129, 10).Value = 5954821.5
$ws.Cells.Item(129, 11).Value = 5369.3334
$ws.Cells.Item(129, 12).Value = 17864464.5
$ws.Cells.Item(129, 13).Value = -369.3334000000004
$ws.Cells.Item(129, 14).Value = -17874464.5
$ws.Cells.Item(131, 8).Value = 908.5540999999999
$ws.Cells.Item(131, 10).Value = 910.0411
$ws.Cells.Item(131, 12).Value = 2730.1233
$ws.Cells.Item(131, 14).Value = -12810.1233
$ws.Cells.Item(132, 8).Value = 2375.7273
$ws.Cells.Item(132, 9).Value = 1585.1111
$ws.Cells.Item(132, 10).Value = 2923.077
$ws.Cells.Item(132, 11).Value = 14265.9999
$ws.Cells.Item(132, 12).Value = 26307.693
$ws.Cells.Item(132, 13).Value = -11735.9999
$ws.Cells.Item(132, 14).Value = -31367.693
$ws.Cells.Item(137, 8).Value = 2288.7646
$ws.Cells.Item(137, 9).Value = 911.2857
$ws.Cells.Item(137, 10).Value = 3253
$ws.Cells.Item(137, 11).Value = 2733.8571
$ws.Cells.Item(137, 12).Value = 9759
$ws.Cells.Item(137, 13).Value = 2366.1429
$ws.Cells.Item(137, 14).Value = -19959

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 65.117645
$ws.Cells.Item(2, 9).Value = 36.77778
$ws.Cells.Item(2, 10).Value = 97
$ws.Cells.Item(2, 11).Value = 36.77778
$ws.Cells.Item(2, 12).Value = 97
$ws.Cells.Item(2, 13).Value = 76.22221999999999
$ws.Cells.Item(2, 14).Value = -323
$ws.Cells.Item(80, 8).Value = 3122.5667
$ws.Cells.Item(80, 9).Value = 2426.7646
$ws.Cells.Item(80, 10).Value = 4032.4614
$ws.Cells.Item(80, 11).Value = 2426.7646
$ws.Cells.Item(80, 12).Value = 4032.4614
$ws.Cells.Item(80, 13).Value = -1428.7646
$ws.Cells.Item(80, 14).Value = -6028.4614
$ws.Cells.Item(83, 8).Value = 3122.5667
$ws.Cells.Item(83, 9).Value = 2426.7646
$ws.Cells.Item(83, 10).Value = 4032.4614
$ws.Cells.Item(83, 11).Value = 12133.823
$ws.Cells.Item(83, 12).Value = 20162.307
$ws.Cells.Item(83, 13).Value = -7141.823
$ws.Cells.Item(83, 14).Value = -30146.307
$ws.Cells.Item(107, 8).Value = 528.6429000000001
$ws.Cells.Item(107, 9).Value = 290.45456
$ws.Cells.Item(107, 10).Value = 1402
$ws.Cells.Item(107, 11).Value = 290.45456
$ws.Cells.Item(107, 12).Value = 1402
$ws.Cells.Item(107, 13).Value = 1629.54544
$ws.Cells.Item(107, 14).Value = -5242
$ws.Cells.Item(132, 8).Value = 50366.805
$ws.Cells.Item(132, 9).Value = 33743
$ws.Cells.Item(132, 10).Value = 101900.6
$ws.Cells.Item(132, 11).Value = 101229
$ws.Cells.Item(132, 12).Value = 305701.8
$ws.Cells.Item(132, 13).Value = -98699
$ws.Cells.Item(132, 14).Value = -310761.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8066941.5
$ws.Cells.Item(7, 9).Value = 10418940
$ws.Cells.Item(7, 10).Value = 2948.4285
$ws.Cells.Item(7, 11).Value = 10418940
$ws.Cells.Item(7, 12).Value = 2948.4285
$ws.Cells.Item(7, 13).Value = -10418828
$ws.Cells.Item(7, 14).Value = -3172.4285
$ws.Cells.Item(82, 8).Value = 1307.375
$ws.Cells.Item(82, 9).Value = 1111.8
$ws.Cells.Item(82, 11).Value = 1111.8
$ws.Cells.Item(82, 13).Value = -750.8
$ws.Cells.Item(85, 8).Value = 1307.375
$ws.Cells.Item(85, 9).Value = 1111.8
$ws.Cells.Item(85, 11).Value = 1111.8
$ws.Cells.Item(85, 13).Value = 136.2
$ws.Cells.Item(93, 8).Value = 1332.1666
$ws.Cells.Item(93, 9).Value = 1098.6
$ws.Cells.Item(93, 11).Value = 1098.6
$ws.Cells.Item(93, 13).Value = 149.4000000000001
$ws.Cells.Item(100, 8).Value = 1848.7858
$ws.Cells.Item(100, 9).Value = 1767.1666
$ws.Cells.Item(100, 10).Value = 1910
$ws.Cells.Item(100, 11).Value = 1767.1666
$ws.Cells.Item(100, 12).Value = 1910
$ws.Cells.Item(100, 13).Value = -1226.1666
$ws.Cells.Item(100, 14).Value = -2992
$ws.Cells.Item(126, 8).Value = 8066941.5
$ws.Cells.Item(126, 9).Value = 10418940
$ws.Cells.Item(126, 10).Value = 2948.4285
$ws.Cells.Item(126, 11).Value = 31256820
$ws.Cells.Item(126, 12).Value = 8845.2855
$ws.Cells.Item(126, 13).Value = -31254350
$ws.Cells.Item(126, 14).Value = -13785.2855
$ws.Cells.Item(132, 8).Value = 85439.086
$ws.Cells.Item(132, 9).Value = 1077.6666
$ws.Cells.Item(132, 10).Value = 169800.5
$ws.Cells.Item(132, 11).Value = 3232.9998
$ws.Cells.Item(132, 12).Value = 509401.5
$ws.Cells.Item(132, 13).Value = -702.9998000000001
$ws.Cells.Item(132, 14).Value = -514461.5
$ws.Cells.Item(136, 8).Value = 54452.59
$ws.Cells.Item(136, 9).Value = 38915.93
$ws.Cells.Item(136, 10).Value = 94000.45
$ws.Cells.Item(136, 11).Value = 116747.79
$ws.Cells.Item(136, 12).Value = 282001.35
$ws.Cells.Item(136, 13).Value = -114197.79
$ws.Cells.Item(136, 14).Value = -287101.35

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2041.1052
$ws.Cells.Item(81, 9).Value = 754
$ws.Cells.Item(81, 11).Value = 1508
$ws.Cells.Item(81, 13).Value = -447
$ws.Cells.Item(84, 8).Value = 2041.1052
$ws.Cells.Item(84, 9).Value = 754
$ws.Cells.Item(84, 11).Value = 7540
$ws.Cells.Item(84, 13).Value = -2236
$ws.Cells.Item(92, 8).Value = 40175
$ws.Cells.Item(92, 10).Value = 40175
$ws.Cells.Item(92, 12).Value = 40175
$ws.Cells.Item(92, 14).Value = -45167
$ws.Cells.Item(136, 8).Value = 44606.766
$ws.Cells.Item(136, 9).Value = 32189.75
$ws.Cells.Item(136, 10).Value = 71096.39999999999
$ws.Cells.Item(136, 11).Value = 96569.25
$ws.Cells.Item(136, 12).Value = 213289.2
$ws.Cells.Item(136, 13).Value = -94019.25
$ws.Cells.Item(136, 14).Value = -218389.2
